$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'244.35"
$ws.Range("E2").Value = "'-1.17%"
$ws.Range("E3").Value = "'3.91%"
$ws.Range("D4").Value = "'5.056"
$ws.Range("E4").Value = "'-0.47%"
$ws.Range("D5").Value = "'0.05683"
$ws.Range("E5").Value = "'1.18%"
$ws.Range("D6").Value = "'6.476"
$ws.Range("E6").Value = "'-0.61%"
$ws.Range("D7").Value = "'0.8205"
$ws.Range("E7").Value = "'0.80%"
$ws.Range("D8").Value = "'0.8393"
$ws.Range("E8").Value = "'-1.08%"
$ws.Range("D9").Value = "'0.1328"
$ws.Range("E9").Value = "'-1.38%"
$ws.Range("D10").Value = "'0.06930"
$ws.Range("E10").Value = "'-0.57%"
$ws.Range("D11").Value = "'0.02856"
$ws.Range("E11").Value = "'-0.66%"
$ws.Range("D12").Value = "'0.09395"
$ws.Range("E12").Value = "'-0.04%"
$ws.Range("D13").Value = "'0.001532"
$ws.Range("E13").Value = "'1.22%"
$ws.Range("D14").Value = "'0.04100"
$ws.Range("D15").Value = "'0.01001"
$ws.Range("E15").Value = "'1.21%"
$ws.Range("D16").Value = "'0.006156"
$ws.Range("E16").Value = "'0.12%"
$ws.Range("D17").Value = "'3.508"
$ws.Range("E17").Value = "'-2.24%"
$ws.Range("D18").Value = "'3.002"
$ws.Range("E18").Value = "'-1.82%"
$ws.Range("D19").Value = "'2.312"
$ws.Range("E19").Value = "'9.13%"
$ws.Range("D21").Value = "'0.03195"
$ws.Range("E21").Value = "'-0.59%"
$ws.Range("D22").Value = "'0.1255"
$ws.Range("E22").Value = "'-4.91%"
$ws.Range("D23").Value = "'3.564"
$ws.Range("E23").Value = "'-5.29%"
$ws.Range("E24").Value = "'1.75%"
$ws.Range("E25").Value = "'-2.59%"
$ws.Range("D26").Value = "'0.003949"
$ws.Range("E26").Value = "'-14.36%"
$ws.Range("D27").Value = "'0.00009799"
$ws.Range("E27").Value = "'2.07%"
$ws.Range("D40").Value = "'0.03731"
$ws.Range("E40").Value = "'1.57%"
$ws.Range("D41").Value = "'0.006099"
$ws.Range("E41").Value = "'-0.82%"
$ws.Range("D42").Value = "'0.1053"
$ws.Range("E42").Value = "'-0.67%"
$ws.Range("D43").Value = "'0.002300"
$ws.Range("E43").Value = "'14.97%"
$ws.Range("D44").Value = "'0.009591"
$ws.Range("E44").Value = "'10.58%"
$ws.Range("D45").Value = "'0.00005215"
$ws.Range("E45").Value = "'-1.47%"
$ws.Range("E46").Value = "'-0.02%"
$ws.Range("E48").Value = "'12.90%"
$ws.Range("E49").Value = "'-0.02%"
$ws.Range("E50").Value = "'-0.02%"
